# Fruta / hortaliza, semanal
# Insert two new weekly observation rows (Frutilla - Vega Modelo de Temuco)
# above the existing row 136, pushing the rest of the table down by 2 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 136-137; everything from the old row 136 onward
# (previously ending at row 163) shifts down to rows 138-165.
$ws.Range("A136:T137").Insert()

# --- New row 136 ---
$ws.Cells.Item(136, 1).Value = 10
$ws.Cells.Item(136, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(136, 3).Value = "La Araucanía"
$ws.Cells.Item(136, 4).Value = 44504
$ws.Cells.Item(136, 5).Value = 9
$ws.Cells.Item(136, 6).Value = "Fruta"
$ws.Cells.Item(136, 7).Value = 100101
$ws.Cells.Item(136, 8).Value = "Berries"
$ws.Cells.Item(136, 9).Value = 100112025
$ws.Cells.Item(136, 10).Value = "Frutilla"
$ws.Cells.Item(136, 11).Value = "Sin especificar"
$ws.Cells.Item(136, 12).Value = "Primera"
$ws.Cells.Item(136, 13).Value = 4000
$ws.Cells.Item(136, 14).Value = 8000
$ws.Cells.Item(136, 15).Value = 9000
$ws.Cells.Item(136, 16).Value = 8550
$ws.Cells.Item(136, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(136, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(136, 19).Value = 1221
$ws.Cells.Item(136, 20).Value = 7

# --- New row 137 ---
$ws.Cells.Item(137, 1).Value = 10
$ws.Cells.Item(137, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(137, 3).Value = "La Araucanía"
$ws.Cells.Item(137, 4).Value = 44504
$ws.Cells.Item(137, 5).Value = 9
$ws.Cells.Item(137, 6).Value = "Fruta"
$ws.Cells.Item(137, 7).Value = 100101
$ws.Cells.Item(137, 8).Value = "Berries"
$ws.Cells.Item(137, 9).Value = 100112025
$ws.Cells.Item(137, 10).Value = "Frutilla"
$ws.Cells.Item(137, 11).Value = "Sin especificar"
$ws.Cells.Item(137, 12).Value = "Segunda"
$ws.Cells.Item(137, 13).Value = 500
$ws.Cells.Item(137, 14).Value = 7000
$ws.Cells.Item(137, 15).Value = 7000
$ws.Cells.Item(137, 16).Value = 7000
$ws.Cells.Item(137, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(137, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(137, 19).Value = 1000
$ws.Cells.Item(137, 20).Value = 7
